$d = $word.ActiveDocument

$d.Content.Find.Execute("623×9=5607", $true, $false, $false, $false, $false, $true, 1, $false, "498×5=2490", 2) | Out-Null
$d.Content.Find.Execute("759×2=1518", $true, $false, $false, $false, $false, $true, 1, $false, "135×8=1080", 2) | Out-Null
$d.Content.Find.Execute("593×2=1186", $true, $false, $false, $false, $false, $true, 1, $false, "890×7=6230", 2) | Out-Null
$d.Content.Find.Execute("542×7=3794", $true, $false, $false, $false, $false, $true, 1, $false, "418×9=3762", 2) | Out-Null
$d.Content.Find.Execute("772×9=6948", $true, $false, $false, $false, $false, $true, 1, $false, "385×9=3465", 2) | Out-Null
$d.Content.Find.Execute("747×9=6723", $true, $false, $false, $false, $false, $true, 1, $false, "373×7=2611", 2) | Out-Null
$d.Content.Find.Execute("816×9=7344", $true, $false, $false, $false, $false, $true, 1, $false, "685×4=2740", 2) | Out-Null
$d.Content.Find.Execute("909×7=6363", $true, $false, $false, $false, $false, $true, 1, $false, "238×9=2142", 2) | Out-Null
$d.Content.Find.Execute("121×8=968", $true, $false, $false, $false, $false, $true, 1, $false, "987×5=4935", 2) | Out-Null
$d.Content.Find.Execute("810×3=2430", $true, $false, $false, $false, $false, $true, 1, $false, "491×5=2455", 2) | Out-Null
$d.Content.Find.Execute("902×6=5412", $true, $false, $false, $false, $false, $true, 1, $false, "294×6=1764", 2) | Out-Null
$d.Content.Find.Execute("636×4=2544", $true, $false, $false, $false, $false, $true, 1, $false, "334×4=1336", 2) | Out-Null
$d.Content.Find.Execute("940×3=2820", $true, $false, $false, $false, $false, $true, 1, $false, "297×2=594", 2) | Out-Null
$d.Content.Find.Execute("259×8=2072", $true, $false, $false, $false, $false, $true, 1, $false, "545×5=2725", 2) | Out-Null
$d.Content.Find.Execute("231×3=693", $true, $false, $false, $false, $false, $true, 1, $false, "109×3=327", 2) | Out-Null
$d.Content.Find.Execute("304×2=608", $true, $false, $false, $false, $false, $true, 1, $false, "287×7=2009", 2) | Out-Null
$d.Content.Find.Execute("430×5=2150", $true, $false, $false, $false, $false, $true, 1, $false, "210×3=630", 2) | Out-Null
$d.Content.Find.Execute("779×7=5453", $true, $false, $false, $false, $false, $true, 1, $false, "892×9=8028", 2) | Out-Null
$d.Content.Find.Execute("288×2=576", $true, $false, $false, $false, $false, $true, 1, $false, "598×2=1196", 2) | Out-Null
$d.Content.Find.Execute("983×7=6881", $true, $false, $false, $false, $false, $true, 1, $false, "413×9=3717", 2) | Out-Null
$d.Content.Find.Execute("938×5=4690", $true, $false, $false, $false, $false, $true, 1, $false, "895×2=1790", 2) | Out-Null
$d.Content.Find.Execute("883×3=2649", $true, $false, $false, $false, $false, $true, 1, $false, "309×9=2781", 2) | Out-Null
$d.Content.Find.Execute("494×2=988", $true, $false, $false, $false, $false, $true, 1, $false, "365×6=2190", 2) | Out-Null
$d.Content.Find.Execute("316×8=2528", $true, $false, $false, $false, $false, $true, 1, $false, "544×8=4352", 2) | Out-Null
$d.Content.Find.Execute("820×5=4100", $true, $false, $false, $false, $false, $true, 1, $false, "508×6=3048", 2) | Out-Null
